# Rename the two "Property" sheets to "DataNode" sheets -- this is the
# core of the commit ("unify the conception of DataNode, DataTable,
# Entity."): Property1 -> DataNode_1, Property2 -> DataNode_2.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$ws1.Name = "DataNode_1"
$ws2.Name = "DataNode_2"

# The author left the workbook with the second sheet (DataNode_2) as the
# active / selected tab.
$ws2.Activate()

# Row-height touch-ups that came along with the resave (header row and the
# wrapped description/header row on each sheet).
$ws1.Rows.Item(1).RowHeight = 27
$ws1.Rows.Item(8).RowHeight = 40.5

$ws2.Rows.Item(1).RowHeight = 27
$ws2.Rows.Item(8).RowHeight = 67.5
